# Rename the "arbitrage" header (columns L:O) to "loser_minus_winner" and add
# a new "winner_minus_loser" block (columns P:S) on both the PRO and CON
# sheets. The new block mirrors the existing "arbitrage"/loser_minus_winner
# block: acar and t_stat are negated (winner - loser = -(loser - winner)),
# while p_value and the significance symbol are identical.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("PRO", "CON")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- header row 1: rename "arbitrage" -> "loser_minus_winner" ---------
    $ws.Range("L1").Value = "loser_minus_winner"

    # new merged header P1:S1 "winner_minus_loser", styled like the other
    # row-1 group headers (bold, thin border, centered/top aligned)
    $ws.Range("P1:S1").Merge()
    $ws.Range("P1").Value = "winner_minus_loser"
    $hdr1 = $ws.Range("P1:S1")
    $hdr1.Font.Bold = $true
    $hdr1.HorizontalAlignment = -4108
    $hdr1.VerticalAlignment = -4160
    $hdr1.Borders.LineStyle = 1

    # --- header row 2 (column labels) -------------------------------------
    $ws.Range("P2").Value = "acar"
    $ws.Range("Q2").Value = "t_stat"
    $ws.Range("R2").Value = "p_value"
    $ws.Range("S2").Value = "symbol"
    $hdr2 = $ws.Range("P2:S2")
    $hdr2.Font.Bold = $true
    $hdr2.HorizontalAlignment = -4108
    $hdr2.VerticalAlignment = -4160
    $hdr2.Borders.LineStyle = 1

    # --- data rows 4-23 -----------------------------------------------
    for ($r = 4; $r -le 23; $r++) {
        $acar = $ws.Cells.Item($r, 12).Value()
        $ws.Cells.Item($r, 16).Value = (0 - $acar)

        $tstat = $ws.Cells.Item($r, 13).Value()
        $ws.Cells.Item($r, 17).Value = (0 - $tstat)

        $pvalue = $ws.Cells.Item($r, 14).Value()
        $ws.Cells.Item($r, 18).Value = $pvalue

        $symbol = $ws.Cells.Item($r, 15).Value()
        $ws.Cells.Item($r, 19).Value = $symbol
    }
}
